# edit.ps1 - apply the 2025-10-17 03:13:58 UTC daily-scrape update to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Data rows -----------------------------------------------------
# Column A holds opportunity IDs that look like plain numbers ("1328588", ...).
# Force text storage first so they keep round-tripping as text cells (as in
# the source data) instead of silently becoming numeric cells.
$ws.Range("A2:A13").NumberFormat = "@"

# Row 2: opportunity 1328588
$ws.Range("A2").Value = "1328588"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328588"
$ws.Range("C2").Value = "My Way​ Operations & Innovation Coordinator​ 2026-2027"
$ws.Range("D2").Value = "40 Düsseldorf, Germany"
$ws.Range("F2").Value = "6 applicants"
$ws.Range("G2").Value = "6 - 18 Months"
$ws.Range("H2").Value = "PwC Global Partnership"

# Row 3: opportunity 1328540
$ws.Range("A3").Value = "1328540"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328540"
$ws.Range("C3").Value = "Media coverage"
$ws.Range("D3").Value = "Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "TAR - Company"

# Row 4: opportunity 1328514
$ws.Range("A4").Value = "1328514"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1328514"
$ws.Range("C4").Value = "AI Intern"
$ws.Range("D4").Value = "Sahibzada Ajit Singh Nagar, Punjab, India"
$ws.Range("F4").Value = "1 applicant"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "CGC JHANJERI MOHALI"

# Row 5: opportunity 1328510
$ws.Range("A5").Value = "1328510"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1328510"
$ws.Range("C5").Value = "Software Developer Intern"
$ws.Range("D5").Value = "Sahibzada Ajit Singh Nagar, Punjab, India"
$ws.Range("F5").Value = "1 applicant"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "CGC JHANJERI MOHALI"

# Row 6: opportunity 1326706
$ws.Range("A6").Value = "1326706"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1326706"
$ws.Range("C6").Value = "Sales Customer Service Support"
$ws.Range("D6").Value = "İstanbul, Türkiye"
$ws.Range("F6").Value = "134 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Aytek Soğutma"

# Row 7: opportunity 1326116
$ws.Range("A7").Value = "1326116"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1326116"
$ws.Range("C7").Value = "Sales and Customer Service Support"
$ws.Range("D7").Value = "İstanbul, Türkiye"
$ws.Range("F7").Value = "106 applicants"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "Tornado Makine Otomotiv İnşaat Sanayi ve Ticaret"

# Row 8: opportunity 1325404
$ws.Range("A8").Value = "1325404"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1325404"
$ws.Range("C8").Value = "Machine Learning Intern"
$ws.Range("D8").Value = "Sahibzada Ajit Singh Nagar, Punjab, India"
$ws.Range("F8").Value = "15 applicants"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "CGC JHANJERI MOHALI"

# Row 9: opportunity 1325403
$ws.Range("A9").Value = "1325403"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1325403"
$ws.Range("C9").Value = "Electrical Engineering Intern"
$ws.Range("D9").Value = "Sahibzada Ajit Singh Nagar, Punjab, India"
$ws.Range("F9").Value = "0 applicants"
$ws.Range("G9").Value = "9 - 12 Weeks"
$ws.Range("H9").Value = "CGC JHANJERI MOHALI"

# Row 10: opportunity 1325318
$ws.Range("A10").Value = "1325318"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1325318"
$ws.Range("C10").Value = "Social Media Planner & Content Creator"
$ws.Range("D10").Value = "القاهرة، محافظة القاهرة‬، مصر"
$ws.Range("F10").Value = "18 applicants"
$ws.Range("G10").Value = "9 - 12 Weeks"
$ws.Range("H10").Value = "The Circle Care"

# Row 11: opportunity 1324728
$ws.Range("A11").Value = "1324728"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1324728"
$ws.Range("C11").Value = "Engineering"
$ws.Range("D11").Value = "Kemalpaşa, İzmir, Türkiye"
$ws.Range("F11").Value = "73 applicants"
$ws.Range("G11").Value = "3 - 6 Months"
$ws.Range("H11").Value = "KEBA OTOMOTİV YEDEKPARÇA SANAYİ VE TİCARİ LİMİTED ŞİRKETİ"

# Row 12: opportunity 1322886
$ws.Range("A12").Value = "1322886"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1322886"
$ws.Range("C12").Value = "Customer Representative"
$ws.Range("D12").Value = "İstanbul, Türkiye"
$ws.Range("F12").Value = "106 applicants"
$ws.Range("G12").Value = "6 - 18 Months"
$ws.Range("H12").Value = "İME HAYAT SAĞLIK TURİZM ANONİM ŞİRKETİ"

# Row 13: opportunity 1317126
$ws.Range("A13").Value = "1317126"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1317126"
$ws.Range("C13").Value = "Social Media Manager"
$ws.Range("D13").Value = "Ahangama, Sri Lanka"
$ws.Range("F13").Value = "18 applicants"
$ws.Range("G13").Value = "9 - 12 Weeks"
$ws.Range("H13").Value = "Surfing Wombats"

# Column A: drop the leftover text-format style so the cells stay plain
# (no style index), matching the source which uses unstyled inline strings.
$ws.Range("A2:A13").Style = "Normal"

# --- 2. PREMIUM highlight moves from row 12 to row 2 -------------------
# Row 2 becomes the premium ("Yes") listing and must pick up the yellow
# highlight style; row 12 reverts to a plain ("No") cell.
$ws.Range("E12").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E2").Value = "Yes"
$ws.Range("E12").Value = "No"
$ws.Range("E12").ClearFormats()

# --- 3. Column width tweaks ---------------------------------------------
$ws.Columns("C").ColumnWidth = 56.166666666666664
$ws.Columns("D").ColumnWidth = 59.166666666666664
$ws.Columns("H").ColumnWidth = 59.166666666666664

# --- 4. Drop the now-stale rows 14-20 (sheet shrinks to A1:H13) --------
$ws.Rows("14:20").Delete()

